# Applies the cryptos list price/volume refresh described in the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.833.24'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '2.926.06'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.88'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.40'
$ws.Range("E6").Value = '  -1.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.631'
$ws.Range("E9").Value = '  +1.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.38'
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.70'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = '3.380.28'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '2.911.84'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.978'
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").Value = '51.810.24'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.29'
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  -2.86%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.75'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.39'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("E26").Value = '  +10.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.10'
$ws.Range("E27").Value = '  +2.27%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +12.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.106'
$ws.Range("E30").Value = '  +12.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.55'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38.54'
$ws.Range("E32").Value = '  -0.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.04'
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '52.18'
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0441'
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  -14.97%  '
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.41'
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.76'
$ws.Range("E41").Value = '  +4.92%  '
$ws.Range("E42").Value = '  +2.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.90'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '120.75'
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("E47").Value = '  -4.15%  '
$ws.Range("D48").Value = '2.139.45'
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.250'
$ws.Range("E49").Value = '  -6.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0335'
$ws.Range("E50").Value = '  +4.69%  '
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.910'
$ws.Range("E51").Value = '  -4.18%  '
